$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A261").Value = "2023-12-13 08:27:56"
$ws.Range("B261").Value = 0.002

$ws.Range("A262").Value = "2023-12-13 08:28:40"
$ws.Range("B262").Value = 0.0032

$ws.Range("A263").Value = "2023-12-13 08:28:43"
$ws.Range("B263").Value = 0.0004
